$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows right after the header row (row 1), shifting the
# existing data (previously rows 2-21) down to rows 12-31.
$insertRange = $ws.Range("A2:C11")
$insertRange.EntireRow.Insert()
# Excel's row insert copies formatting from the row above (the bold
# header row); clear it so the new data rows stay unformatted, matching
# the rest of the data rows.
$insertRange.ClearFormats()

# New data for the inserted rows (A2:C11)
$newData = @(
    @(-0.0058032199740409, 0.0117591563612222, 0.0122173046693205),
    @(-0.0114537235349416, -0.0096211275085806, -0.0482583530247211),
    @(0.0445931628346443, 0.1122464910149574, -0.0378736443817615),
    @(0.0612392425537109, 0.09758572280406951, -0.0021380283869802),
    @(0.0088575463742017, 0.1237002089619636, 0.0548251569271087),
    @(-0.0221438650041818, 0.0061086523346602, 0.0325285755097866),
    @(0.0332921557128429, -0.0615446716547012, 0.093156948685646),
    @(-0.4489859640598297, -1.353219270706177, 0.4497495293617248),
    @(-0.3888157308101654, -3.63083028793335, -0.1369865238666534),
    @(-0.6565274000167847, -2.371837139129639, 0.1600466966629028)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
